# "Drop in RMI script files" - applies the Discount Rate workbook edits:
#  1. Remove the "Texas Notes" worksheet (it was specific to the old Texas
#     scenario and isn't used any more).
#  2. Update the discount rate value on the "DR" sheet from 0.0587 to 0.03.
#  3. Update the selections that were left on the sheets when the file was
#     last saved (About -> A16:A18, DR -> B1).

$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# 1. Drop the "Texas Notes" worksheet entirely.
$wsNotes = $wb.Worksheets.Item("Texas Notes")
$wsNotes.Delete()

# 2. Update the DR discount rate value (B2) to 3%, and fix up its selection.
$wsDR = $wb.Worksheets.Item("DR")
$wsDR.Range("B2").Value = 0.03
$wsDR.Activate()
$wsDR.Range("B1").Select()

# 3. Restore "About" as the selected/active tab, with its new selection.
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate()
$wsAbout.Range("A16:A18").Select()
